$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row before the current last table row (row 27) ---
# This pushes the old "closing" row (27, with the bottom-border style) down
# to row 28, and also pushes the footer rows (old 32/33) down to 33/34,
# matching the diff (one new worker period row was added for Keiner).
$ws.Rows.Item(27).Insert()

# Copy the formatting of the row above (26, a "normal" data row) into the
# newly inserted row 27 so it matches the rest of the table (borders etc.)
$ws.Range("B26:J26").Copy()
$ws.Range("B27:J27").PasteSpecial(-4122)

# --- Yamile Reyes Bayona (rows 16-18): periods now listed most-recent-first ---
$ws.Range("E16").Value = "2310"
$ws.Range("E17").Value = "2309"
$ws.Range("E18").Value = "2308"

# --- Keiner Stiven Aguilar Jimenez (rows 19-28): new period 2507 added,
#     whole block now listed most-recent-first ---
$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1065871889"
$ws.Range("D27").Value = "KEINER STIVEN AGUILAR JIMENEZ"
$ws.Range("F27").Value = 52000
$ws.Range("G27").Value = 1300000

$ws.Range("E19").Value = "2507"
$ws.Range("E20").Value = "2506"
$ws.Range("E21").Value = "2505"
$ws.Range("E22").Value = "2504"
$ws.Range("E23").Value = "2503"
$ws.Range("E24").Value = "2502"
$ws.Range("E25").Value = "2501"
$ws.Range("E26").Value = "2412"
$ws.Range("E27").Value = "2411"
$ws.Range("E28").Value = "2410"

# --- Update the summary fields above the table ---
# Valor Mora total (3 x 46400 + 10 x 52000)
$ws.Range("E11").Value = 659200
# Cant. Periodos (now 13 total period rows instead of 12)
$ws.Range("F13").Value = 13
